$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15
$ws.Range("D15").Value = 44806
$ws.Range("H15").Value = "Madrigal"
$ws.Range("I15").Value = "Segunda"
$ws.Range("J15").Value = 80
$ws.Range("K15").Value = 13000
$ws.Range("L15").Value = 13000
$ws.Range("M15").Value = 13000
$ws.Range("N15").Value = "$/caja 50 unidades"
$ws.Range("P15").Value = 260
$ws.Range("Q15").Value = 50

# Row 16
$ws.Range("D16").Value = 44803
$ws.Range("H16").Value = "Madrigal"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 100
$ws.Range("K16").Value = 14000
$ws.Range("L16").Value = 15000
$ws.Range("M16").Value = 14500
$ws.Range("N16").Value = "$/caja 40 unidades"
$ws.Range("P16").Value = 362
$ws.Range("Q16").Value = 40

# Row 17
$ws.Range("D17").Value = 44495
$ws.Range("H17").Value = "Madrigal"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 130
$ws.Range("K17").Value = 11000
$ws.Range("L17").Value = 11000
$ws.Range("M17").Value = 11000
$ws.Range("N17").Value = "$/caja 40 unidades"
$ws.Range("P17").Value = 275
$ws.Range("Q17").Value = 40

# Row 18
$ws.Range("D18").Value = 44810
$ws.Range("H18").Value = "Madrigal"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 100
$ws.Range("K18").Value = 13500
$ws.Range("L18").Value = 14000
$ws.Range("M18").Value = 13750
$ws.Range("N18").Value = "$/caja 40 unidades"
$ws.Range("P18").Value = 344
$ws.Range("Q18").Value = 40

# Row 19
$ws.Range("D19").Value = 44845
$ws.Range("H19").Value = "Española"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 120
$ws.Range("K19").Value = 10000
$ws.Range("L19").Value = 10000
$ws.Range("M19").Value = 10000
$ws.Range("N19").Value = "$/caja 30 unidades"
$ws.Range("P19").Value = 333
$ws.Range("Q19").Value = 30
